$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $result = $d.Content.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $result) {
        Write-Host "WARNING: Find/Replace failed for: $findText"
    }
    return $result
}

# 1. " web crawlers, " -> " web crawlers " (drop comma before "followed")
Replace-Text "web crawlers, followed" "web crawlers followed"

# 2. "together by making it a web application" -> "together by making a web application"
Replace-Text "together by making it a web application" "together by making a web application"

# 3. "each of them is able" -> "each client is able"
Replace-Text "each of them is able" "each client is able"

# 4. Rework the tail of the "web crawlers" paragraph plus the two following paragraphs in
#    one pass: shorten the first paragraph, revise the "Next, I went on..." paragraph, and
#    revise the "Currently for my progress..." paragraph (including its trailing GitHub text).
$findText = "their crawlers to crawl whichever forums/websites they want, returning the specific desired data in an organized manner back to the client.`rNext, I went on to design and test out the overall architecture for weeks and it was one of the challenge I had because it was complicated. Finally, after numerous attempts, I managed to come out with an architecture integrating Django, Scrapy and MongoDB together in a multi-threaded and also a multi-process environment.    `rCurrently for my progress, I have a planned architecture and also have already set up the necessary environment for the development of the crawler with a basic working prototype which is able to handle multiple client connections with multiple crawler processes. You can view the source code and my commits at this GitHub link: "
$replaceText = "their own crawlers.^pNext, I went on to design and test out the overall architecture for several weeks and this was one of the main challenges I had because it was complicated but finally after numerous attempts, I managed to come out with an architecture integrating Django, Scrapy and MongoDB altogether in a multi-threaded and also multi-process environment.    ^pCurrently for my progress, I have developed an architecture and also set up the necessary environment for the development of the crawler with a basic working prototype which is able to handle multiple client connections with multiple crawler processes. You can view the source code and my commits at my GitHub repository link: "
Replace-Text $findText $replaceText

# 5. "some of the more complex" -> "some of the complex"
Replace-Text "some of the more complex" "some of the complex"

# 6. "using older versions" -> "using old versions"
Replace-Text "using older versions" "using old versions"

# 7. "which has already since changed drastically." -> "which has already changed drastically."
Replace-Text "which has already since changed drastically." "which has already changed drastically."

# 8. "to operate, but to integrate" -> "to operate but to integrate"
Replace-Text "to operate, but to integrate" "to operate but to integrate"

# 9. "Next Step" heading -> "Next Steps"
Replace-Text "Next Step" "Next Steps"

# 10. "Controlling on the number of crawlers and CPU resource for each user"
#     -> "Controlling the number of crawlers and CPU resource for each client"
Replace-Text "Controlling on the number of crawlers and CPU resource for each user" "Controlling the number of crawlers and CPU resource for each client"

# 11. "Auto detection of structured data based on user input"
#     -> "Auto detection of structured data based on client input"
Replace-Text "Auto detection of structured data based on user input" "Auto detection of structured data based on client input"
